# Refactor the "currency_conversions" sheet so that instead of a single
# "foreign_amount" column it has two explicit columns: "source_amount" and
# "target_amount" (plus a matching "target_fees" next to the existing
# "source_fees"). Columns end up as:
#   date | source_amount | source_fees | source_currency |
#   target_amount | target_fees | target_currency | comment

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# Make room for the two new columns ("target_amount" and "target_fees") by
# inserting them right before the existing "target_currency" column (E).
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).Insert()

# Set the new headers. "target_fees" is written first, then "source_amount"
# (renaming the old "foreign_amount" header), then "target_amount" - this
# ordering only affects the order new strings are interned, not the result.
$ws.Range("F1").Value = "target_fees"
$ws.Range("B1").Value = "source_amount"
$ws.Range("E1").Value = "target_amount"

# Fill in the new target_amount / target_fees values for the two existing
# currency-conversion rows.
$ws.Range("E2").Value = -1
$ws.Range("F2").Value = 0
$ws.Range("E3").Value = -1
$ws.Range("F3").Value = 0

# Make "currency_conversions" the active tab (it was previously
# "money_transfers").
$ws.Activate()
